# Generate Report for Handoff
# Updates the localization-status report: marks the "ht" (handoff) priority
# for the rows whose handoff just completed, and refreshes the "Latest
# HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps for those
# same rows on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 12, 14)

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $cell = $overview.Range("G$r")
    $cell.NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $cell.Value = "2016-08-17 10:20:36"
}

# zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "ht"
    $hcell = $zhcn.Range("H$r")
    $hcell.NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $hcell.Value = "2016-08-17 10:20:30"
}

# de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("E$r").Value = "ht"
    $hcell2 = $dede.Range("H$r")
    $hcell2.NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $hcell2.Value = "2016-08-17 10:20:36"
}
